$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34, shifting existing rows 34:182 down to 35:183
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with its data
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44676
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112037
$ws.Range("G34").Value = "Cebollín"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 800
$ws.Range("K34").Value = 1100
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = 1150
$ws.Range("N34").Value = "$/paquete 6 unidades"
$ws.Range("O34").Value = "Provincia del Elquí"
$ws.Range("P34").Value = 192
$ws.Range("Q34").Value = 6
$ws.Range("R34").Value = "Hortaliza"
